$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("effort")

# Carry the date format (style used by column A, e.g. A32) down onto the
# three new rows before filling in values.
$ws.Range("A32").Copy()
$ws.Range("A33:A35").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 33: 30-Oct-2012 - manual continued (draft, reused entry), plus review effort
$ws.Range("A33").Value = 41212
$ws.Range("B33").Value = 2.75
$ws.Range("C33").Value = 2
$ws.Range("D33").Value = "Manual continued"

# Row 34: 31-Oct-2012 - manual review work
$ws.Range("A34").Value = 41213
$ws.Range("B34").Value = 2.75
$ws.Range("D34").Value = "Manual reviewed"

# Row 35: 01-Nov-2012 - manual review work completed
$ws.Range("A35").Value = 41214
$ws.Range("B35").Value = 1.5
$ws.Range("D35").Value = "Manual reviewed"

# Reflect the scrolled/selected state after the edit (matches the author's
# saved view: scrolled down to row 19, final cell D35 selected).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D35").Select()
